# QE_holdings.xlsx update: refresh model-holdings weights / percent-change
# figures and bump the "as of" date in the confidentiality footnote.
#
# The sheet carries protection, so we temporarily unprotect it, apply the
# value + text edits, then restore protection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Footnote text: "as of 2021-07-13" -> "as of 2021-07-14" -------------
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."

# --- Weight (D) / Percent Change (E) figures, rows 2-35 -------------------
$ws.Range("D2").Value = 0.09476064219138863
$ws.Range("E2").Value = 0.02410052183466105
$ws.Range("D3").Value = 0.07948677872459903
$ws.Range("E3").Value = 0.00544522741832143
$ws.Range("D4").Value = 0.05201464527914646
$ws.Range("E4").Value = 0.001174755803076133
$ws.Range("D5").Value = 0.0513409398071106
$ws.Range("E5").Value = -0.003223380251423436
$ws.Range("D6").Value = 0.04799108326832729
$ws.Range("E6").Value = 0.01501959975242428
$ws.Range("D7").Value = 0.04183041943575929
$ws.Range("E7").Value = -0.003405075489881004
$ws.Range("D8").Value = 0.03602379398163046
$ws.Range("E8").Value = 0.007032271490440944
$ws.Range("D9").Value = 0.03830799924468041
$ws.Range("E9").Value = 0.006793879600637931
$ws.Range("D10").Value = 0.03380352213729985
$ws.Range("E10").Value = 0.006899985773225259
$ws.Range("D11").Value = 0.03552038903202399
$ws.Range("E11").Value = -0.009079180006689946
$ws.Range("D12").Value = 0.03486111101288498
$ws.Range("E12").Value = -0.01266721576869545
$ws.Range("D13").Value = 0.03044532030713501
$ws.Range("E13").Value = -0.02508780732563975
$ws.Range("D14").Value = 0.03139173292804449
$ws.Range("E14").Value = 0.006844346317615546
$ws.Range("D15").Value = 0.03256445853557109
$ws.Range("E15").Value = 0.01847749602564464
$ws.Range("D16").Value = 0.03117178499417598
$ws.Range("E16").Value = -0.001252382248843054
$ws.Range("D17").Value = 0.02906070916588813
$ws.Range("E17").Value = 0.007373877491421377
$ws.Range("D18").Value = 0.02836861576368964
$ws.Range("E18").Value = -0.02267627965417174
$ws.Range("D19").Value = 0.02371731810624121
$ws.Range("E19").Value = -0.0008647526807331607
$ws.Range("D20").Value = 0.02062758005645095
$ws.Range("E20").Value = 0.004100552683187653
$ws.Range("D21").Value = 0.02152802285195382
$ws.Range("E21").Value = -0.02217477003942181
$ws.Range("D22").Value = 0.02179054591642969
$ws.Range("E22").Value = 0.02253725917848048
$ws.Range("D23").Value = 0.02091443176698815
$ws.Range("E23").Value = -0.006154387198874511
$ws.Range("D24").Value = 0.01839387087853994
$ws.Range("E24").Value = -0.0003537318712416582
$ws.Range("D25").Value = 0.02216820571284701
$ws.Range("E25").Value = 0.01389686459170791
$ws.Range("D26").Value = 0.02018994731858008
$ws.Range("E26").Value = 0.007566204287515976
$ws.Range("D27").Value = 0.01960888872544062
$ws.Range("E27").Value = 0.01444111027756945
$ws.Range("D28").Value = 0.01846077465520863
$ws.Range("E28").Value = -0.01247366203118416
$ws.Range("D29").Value = 0.02057057747294676
$ws.Range("E29").Value = -0.0003094250881862104
$ws.Range("D30").Value = 0.01145709494749185
$ws.Range("E30").Value = -0.02017283950617288
$ws.Range("D31").Value = 0.008520401054901408
$ws.Range("E31").Value = -0.001460871874896186
$ws.Range("D32").Value = 0.007647681600259127
$ws.Range("E32").Value = 0.01344603092402186
$ws.Range("D33").Value = 0.008557601252076597
$ws.Range("E33").Value = 0.006313945224045936
$ws.Range("D34").Value = 0.006903111874288794
$ws.Range("E34").Value = -0.006515859355790576
$ws.Range("E35").Value = 0.00307630071620979

# Restore worksheet protection (sheet-content protection was present before
# the edit).
$ws.Protect("D382", $true, $true, $true, $false, $true, $false, $false)
